# Fills in the previously-blank "Answer" cells across the four tables of
# the Investments workbook with the computed results, and marks the big
# total-investment figure with an integer ("0") number format.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Table -1.1")
$ws2 = $wb.Worksheets.Item("Table - 2.1")
$ws3 = $wb.Worksheets.Item("Table-3.1")
$ws4 = $wb.Worksheets.Item("Table-5.1")

# --- Table -1.1 : "Understand the Data Set" answers ---
$ws1.Cells.Item(5,3).Value = 66370
$ws1.Cells.Item(6,3).Value = 66368
$ws1.Cells.Item(7,3).Value = 'permalink'
$ws1.Cells.Item(8,3).Value = 'Y'
$ws1.Cells.Item(9,3).Value = 114942

# --- Table - 2.1 : average funding amounts per investment type ---
$ws2.Cells.Item(5,3).Value = '11724223  USD'
$ws2.Cells.Item(6,3).Value = '971575 USD'
$ws2.Cells.Item(7,3).Value = '747793 USD'
$ws2.Cells.Item(8,3).Value = '73938484 USD'
$ws2.Cells.Item(9,3).Value = 'Venture investment'

# --- Table-3.1 : top English-speaking countries ---
$ws3.Cells.Item(5,3).Value = ' United States'
$ws3.Cells.Item(6,3).Value = 'United Kingdom'
$ws3.Cells.Item(7,3).Value = 'India'

# --- Table-5.1 : sector-wise investment analysis ---
$ws4.Cells.Item(5,3).Value = 12092
$ws4.Cells.Item(5,4).Value = 622
$ws4.Cells.Item(5,5).Value = 328

$ws4.Cells.Item(6,3).Value = 108002096957
$ws4.Cells.Item(6,4).Value = 5394078692
$ws4.Cells.Item(6,5).Value = 2949543602
$ws4.Cells.Item(6,3).NumberFormat = "0"

$ws4.Cells.Item(7,3).Value = 'Others'
$ws4.Cells.Item(7,4).Value = 'Others'
$ws4.Cells.Item(7,5).Value = 'Others'

$ws4.Cells.Item(8,3).Value = 'Social, Finance, Analytics, Advertising'
$ws4.Cells.Item(8,4).Value = 'Social, Finance, Analytics, Advertising'
$ws4.Cells.Item(8,5).Value = 'Social, Finance, Analytics, Advertising'

$ws4.Cells.Item(9,3).Value = 'Cleantech / Semiconductors'
$ws4.Cells.Item(9,4).Value = 'Cleantech / Semiconductors'
$ws4.Cells.Item(9,5).Value = 'News, Search and Messaging'

$ws4.Cells.Item(10,3).Value = 2957
$ws4.Cells.Item(10,4).Value = 148
$ws4.Cells.Item(10,5).Value = 110

$ws4.Cells.Item(11,3).Value = 2718
$ws4.Cells.Item(11,4).Value = 133
$ws4.Cells.Item(11,5).Value = 60

$ws4.Cells.Item(12,3).Value = 2355
$ws4.Cells.Item(12,4).Value = 130
$ws4.Cells.Item(12,5).Value = 52

$ws4.Cells.Item(13,3).Value = 'Virtustream'
$ws4.Cells.Item(13,4).Value = 'Electric Cloud'
$ws4.Cells.Item(13,5).Value = 'FirstCry.com'

$ws4.Cells.Item(14,3).Value = 'SST Inc. (Formerly ShotSpotter)'
$ws4.Cells.Item(14,4).Value = 'Celltick Technologies'
$ws4.Cells.Item(14,5).Value = 'Manthan Systems'

# --- restore the sheet/selection state left behind by the edit ---
$ws2.Range("C9").Select()
$ws3.Range("C8").Select()
$ws4.Range("E9").Select()

$ws1.Activate()
$ws1.Range("B5").Select()

Write-Output "done"
